$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 255. This shifts the existing
# rows 255-298 down to 256-299 and leaves an empty row 255 ready to be
# populated with the new weekly data point.
$ws.Rows(255).Insert()

# Populate the new row 255 with the new observation (mirrors the fixed
# columns used throughout the sheet for this market/product).
$ws.Range("A255").Value = 6
$ws.Range("B255").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C255").Value = "Metropolitana"
$ws.Range("D255").Value = 44505
$ws.Range("E255").Value = 13
$ws.Range("F255").Value = 100112043
$ws.Range("G255").Value = "Pepino ensalada"
$ws.Range("H255").Value = "Sin especificar"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 400
$ws.Range("K255").Value = 7000
$ws.Range("L255").Value = 8000
$ws.Range("M255").Value = 7575
$ws.Range("N255").Value = "$/caja 70 unidades"
$ws.Range("O255").Value = "Región de O'Higgins"
$ws.Range("P255").Value = 108
$ws.Range("Q255").Value = 70
$ws.Range("R255").Value = "Hortaliza"

# Match the date-number formatting used by the rest of column D.
$ws.Range("D255").NumberFormat = "YYYY-MM-DD HH:MM:SS"
